$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efnb1", "Ephb1", "ECs", 3, 1, 9.712432666666666, 29.137298, 0.4639063029983291, 0.463906302998329, 3, 1, 2.814284333333333, 8.442853, 0.7633363087964763, 0.7633363087964762, 27.33354709235489, 246.001923831194, 0.3541165249581643, 0.3541165249581641),
    @("ECs", "Efnb1", "Ephb1", "FAPs", 3, 1, 9.712432666666666, 29.137298, 0.4639063029983291, 0.463906302998329, 2, 0.6666666666666666, 0.1814093333333333, 0.5442279999999999, 0.04920481176963387, 0.04920481176963386, 1.761925935104889, 15.857333415944, 0.02282642231777952, 0.02282642231777951),
    @("ECs", "Efnb1", "Ephb1", "sCs", 3, 1, 9.712432666666666, 29.137298, 0.4639063029983291, 0.463906302998329, 3, 1, 0.6911273333333333, 2.073382, 0.1874588794338899, 0.1874588794338899, 6.712527689092888, 60.412749201836, 0.08696335572238538, 0.08696335572238535),
    @("FAPs", "Efnb1", "Ephb1", "ECs", 3, 1, 9.124904999999998, 27.37471499999999, 0.4358435305594535, 0.4358435305594534, 3, 1, 2.814284333333333, 8.442853, 0.7633363087964763, 0.7633363087964762, 25.68007718465499, 231.1206946618949, 0.3326951918300774, 0.3326951918300773),
    @("FAPs", "Efnb1", "Ephb1", "FAPs", 3, 1, 9.124904999999998, 27.37471499999999, 0.4358435305594535, 0.4358435305594534, 2, 0.6666666666666666, 0.1814093333333333, 0.5442279999999999, 0.04920481176963387, 0.04920481176963386, 1.65534293278, 14.89808639501999, 0.02144559888219058, 0.02144559888219057),
    @("FAPs", "Efnb1", "Ephb1", "sCs", 3, 1, 9.124904999999998, 27.37471499999999, 0.4358435305594535, 0.4358435305594534, 3, 1, 0.6911273333333333, 2.073382, 0.1874588794338899, 0.1874588794338899, 6.306471259569999, 56.75824133612999, 0.08170273984718548, 0.08170273984718548),
    @("M1", "Efnb1", "Ephb1", "ECs", 3, 1, 0.318184, 0.9545520000000001, 0.01519779525677573, 0.01519779525677573, 3, 1, 2.814284333333333, 8.442853, 0.7633363087964763, 0.7633363087964762, 0.8954602463173333, 8.059142216856, 0.01160102893315178, 0.01160102893315178),
    @("M1", "Efnb1", "Ephb1", "FAPs", 3, 1, 0.318184, 0.9545520000000001, 0.01519779525677573, 0.01519779525677573, 2, 0.6666666666666666, 0.1814093333333333, 0.5442279999999999, 0.04920481176963387, 0.04920481176963386, 0.05772154731733333, 0.519493925856, 0.0007478046549230844, 0.0007478046549230843),
    @("M1", "Efnb1", "Ephb1", "sCs", 3, 1, 0.318184, 0.9545520000000001, 0.01519779525677573, 0.01519779525677573, 3, 1, 0.6911273333333333, 2.073382, 0.1874588794338899, 0.1874588794338899, 0.2199056594293333, 1.979150934864, 0.002848961668700866, 0.002848961668700866),
    @("M2", "Efnb1", "Ephb1", "ECs", 3, 1, 0.353022, 1.059066, 0.01686180347577968, 0.01686180347577968, 3, 1, 2.814284333333333, 8.442853, 0.7633363087964763, 0.7633363087964762, 0.993504283922, 8.941538555298, 0.01287122682485326, 0.01287122682485325),
    @("M2", "Efnb1", "Ephb1", "FAPs", 3, 1, 0.353022, 1.059066, 0.01686180347577968, 0.01686180347577968, 2, 0.6666666666666666, 0.1814093333333333, 0.5442279999999999, 0.04920481176963387, 0.04920481176963386, 0.06404148567199999, 0.576373371048, 0.0008296818661222976, 0.0008296818661222972),
    @("M2", "Efnb1", "Ephb1", "sCs", 3, 1, 0.353022, 1.059066, 0.01686180347577968, 0.01686180347577968, 3, 1, 0.6911273333333333, 2.073382, 0.1874588794338899, 0.1874588794338899, 0.243983153468, 2.195848381212, 0.003160894784804129, 0.003160894784804129),
    @("sCs", "Efnb1", "Ephb1", "ECs", 3, 1, 1.427651, 4.282953, 0.06819056770966213, 0.06819056770966211, 3, 1, 2.814284333333333, 8.442853, 0.7633363087964763, 0.7633363087964762, 4.017815842767666, 36.160342584909, 0.05205233625022967, 0.05205233625022965),
    @("sCs", "Efnb1", "Ephb1", "FAPs", 3, 1, 1.427651, 4.282953, 0.06819056770966213, 0.06819056770966211, 2, 0.6666666666666666, 0.1814093333333333, 0.5442279999999999, 0.04920481176963387, 0.04920481176963386, 0.2589892161426666, 2.330902945284, 0.003355304048618399, 0.003355304048618397),
    @("sCs", "Efnb1", "Ephb1", "sCs", 3, 1, 1.427651, 4.282953, 0.06819056770966213, 0.06819056770966211, 3, 1, 0.6911273333333333, 2.073382, 0.1874588794338899, 0.1874588794338899, 0.9866886285606666, 8.880197657046, 0.01278292741081406, 0.01278292741081406)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
